# Modified DSL for EB
# - Clear the "Pass"/"Fail" Results column (J2:J15) since the Results
#   column no longer carries a canned Pass/Fail value.
# - Update the package/activity used to relaunch the app after toggling
#   wifi from the old Rhomobile compliance-test package to the new
#   Enterprise Browser package (rows 4, 14, 15 / cells G4, G14, G15).
# - Move the active selection back to G1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2:J15").ClearContents()

$ws.Range("G4").Value = "wait(3);
validate1;
link_Click(signal_test_link);
validate2;
SelectTestToRun(VT200_0852_string);
ClickRunTest(runtest_top_xpath);
validate3;
ClickRunTest(runtest_bottom_xpath);
wait(2);
wifi_Mode(OFF);
wait(2);
wifi_Mode(ON);
wait(2);
press_Key(Home);
launch_App_Device(com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);
wait(2);
validate4;
checkCallbackValues(essid_xpath);
checkCallbackValues(ipaddress_xpath);
checkCallbackValues(signalStrength_xpath);
signalCallbackcount(results_id);

"

$ws.Range("G14").Value = "wait(3);
validate1;
link_Click(signal_test_link);
validate2;
SelectTestToRun(VT200_0862_string);
ClickRunTest(runtest_top_xpath);
validate3;
ClickRunTest(runtest_bottom_xpath);
wifi_Mode(OFF);
wifi_Mode(ON);
press_Key(Home);
launch_App_Device(com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);
validate4;
wait(15);
checkstopwlanStatus(results_id);"

$ws.Range("G15").Value = "wait(3);
validate1;
link_Click(signal_test_link);
validate2;
SelectTestToRun(VT200_0864_string);
ClickRunTest(runtest_top_xpath);
validate3;
wifi_Mode(OFF);
wait(2);
press_Key(Home);
launch_App_Device(com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);
wait(2);
ClickRunTest(runtest_bottom_xpath);
wait(2);
validate4;
checkCallbackValues(ipwlandisable_xpath);
wifi_Mode(ON);
press_Key(Home);"

[void]$ws.Range("G1").Select()
